# Deploying to gh-pages: refresh the "Metadata" summary sheet and the
# root-extension row on the "Elements" sheet for the new IG build.

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (Property / Value table) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# New publication date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher is now filled in (was blank)
$meta.Range("B9").Value = "Alvearie Team"

# The old duplicate "Contact / No display for ContactDetail" rows (10 & 11)
# are replaced by a single "Jurisdiction / United States of America" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

# --- Sheet "Elements" (StructureDefinition element table) ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row (row 2): Short / Definition now reflect the profile
$elements.Range("K2").Value = "Rx Formulary Indicator"
$elements.Range("L2").Value = "Indicates whether the prescription drug is included in the formulary"
